# Auto-generated edit script: updates crypto price/volume table
# to reflect the Sun Jan 14 08:39:05 UTC 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text/string updates (safe from Excel's numeric auto-detection) ---
$ws.Range("D2").Value = '42.931.59'
$ws.Range("E2").Value = '  -0.56%  '
$ws.Range("D3").Value = '2.551.68'
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("E5").Value = '  +1.46%  '
$ws.Range("E6").Value = '  +3.89%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  -1.14%  '
$ws.Range("E10").Value = '  +2.23%  '
$ws.Range("E11").Value = '  +1.17%  '
$ws.Range("E12").Value = '  -0.21%  '
$ws.Range("E13").Value = '  +5.77%  '
$ws.Range("D14").Value = '2.943.77'
$ws.Range("E14").Value = '  -0.25%  '
$ws.Range("D15").Value = '2.493.50'
$ws.Range("E15").Value = '  -2.02%  '
$ws.Range("E16").Value = '  +6.17%  '
$ws.Range("E17").Value = '  -0.39%  '
$ws.Range("D18").Value = '42.965.91'
$ws.Range("E18").Value = '  -0.61%  '
$ws.Range("E19").Value = '  +3.83%  '
$ws.Range("D20").Value = '0.0₃0993'
$ws.Range("E20").Value = '  +1.29%  '
$ws.Range("E21").Value = '  -0.58%  '
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("E23").Value = '  -2.37%  '
$ws.Range("E24").Value = '  +1.03%  '
$ws.Range("E25").Value = '  -2.30%  '
$ws.Range("E26").Value = '  -5.12%  '
$ws.Range("E27").Value = '  -0.16%  '
$ws.Range("E28").Value = '  +0.99%  '
$ws.Range("E29").Value = '  +0.37%  '
$ws.Range("E30").Value = '  -1.04%  '
$ws.Range("E31").Value = '  +1.83%  '
$ws.Range("E32").Value = '  +2.55%  '
$ws.Range("E33").Value = '  -1.04%  '
$ws.Range("E34").Value = '  -1.00%  '
$ws.Range("E35").Value = '  +0.82%  '
$ws.Range("E36").Value = '  +12.07%  '
$ws.Range("E37").Value = '  -2.52%  '
$ws.Range("E38").Value = '  +11.08%  '
$ws.Range("E39").Value = '  -0.17%  '
$ws.Range("E40").Value = '  -0.33%  '
$ws.Range("E41").Value = '  +33.42%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("B43").Value = 'NEARProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("E43").Value = '  -1.15%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("E44").Value = '  -2.23%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '2.084.85'
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("E46").Value = '  -0.01%  '
$ws.Range("E47").Value = '  +0.67%  '
$ws.Range("E48").Value = '  +0.57%  '
$ws.Range("D49").Value = '2.801.09'
$ws.Range("E49").Value = '  -0.17%  '
$ws.Range("E50").Value = '  +7.83%  '
$ws.Range("E51").Value = '  -1.56%  '

# --- Price cells whose new text reads as a plain number (e.g. "37.00", "0.999") ---
# Force Text number format first so Excel keeps the exact string (incl. trailing
# zeros) instead of silently casting it to a float cell.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.32'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.19'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.545'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.00'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0819'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.75'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.03'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.872'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.62'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.59'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.96'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '253.98'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '28.11'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.16'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '37.82'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.11'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '158.49'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.16'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0804'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '19.08'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.92'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.11'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.90'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.42'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0306'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '86.44'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '74.82'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '103.33'

# --- Restore the original (unstyled/default) cell format on those price cells; ---
# --- only the text NumberFormat flag was needed to pin the type during entry.   ---
$fmtRef = $ws.Range("D24")
$fmtRef.Copy()
$targetCells = @("D5","D6","D9","D10","D11","D12","D16","D17","D19","D21","D22","D23","D26","D27","D28","D29","D31","D32","D33","D35","D36","D38","D41","D42","D43","D44","D47","D50","D51")
foreach ($cellRef in $targetCells) {
    $ws.Range($cellRef).PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = $false

